# Mark all slides except slide4 and slide5 as hidden (p:sld/@show="0").
# This mirrors the commit's OOXML diff, which adds show="0" to the
# <p:sld> root element of slides 1,2,3,6,7,8,9,10,11,12,13,14,15,16,17,18.

$p = $ppt.ActivePresentation

$hiddenSlideNumbers = @(1,2,3,6,7,8,9,10,11,12,13,14,15,16,17,18)

foreach ($n in $hiddenSlideNumbers) {
    $s = $p.Slides.Item($n)
    $s.SlideShowTransition.Hidden = 1
}
